# New crime data collected - weekly CompStat update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report Volume/Number and week date range) ---
$ws.Range("A8").Value = "Volume 32   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/29/2025  Through  10/5/2025"

# --- Weekly crime statistics table (rows 14-30): plain numeric value updates ---
$ws.Range("C14").Value = 1
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = -80
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 15
$ws.Range("H14").Value = -73.333333333333
$ws.Range("I14").Value = 82
$ws.Range("J14").Value = 96
$ws.Range("K14").Value = -14.583333333333
$ws.Range("L14").Value = -23.364485981308
$ws.Range("M14").Value = -17.171717171717
$ws.Range("N14").Value = -78.756476683937
$ws.Range("C15").Value = 10
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = 42.857142857142
$ws.Range("F15").Value = 43
$ws.Range("G15").Value = 35
$ws.Range("H15").Value = 22.857142857142
$ws.Range("I15").Value = 409
$ws.Range("J15").Value = 321
$ws.Range("K15").Value = 27.414330218068
$ws.Range("L15").Value = 39.115646258503
$ws.Range("M15").Value = 76.293103448275
$ws.Range("N15").Value = -26.173285198556
$ws.Range("C16").Value = 90
$ws.Range("D16").Value = 100
$ws.Range("E16").Value = -10
$ws.Range("F16").Value = 408
$ws.Range("G16").Value = 400
$ws.Range("H16").Value = 2
$ws.Range("I16").Value = 3723
$ws.Range("J16").Value = 3834
$ws.Range("K16").Value = -2.895148669796
$ws.Range("L16").Value = 0.107555794568
$ws.Range("M16").Value = 10.573210573210
$ws.Range("N16").Value = -69.603200522534
$ws.Range("C17").Value = 177
$ws.Range("D17").Value = 146
$ws.Range("E17").Value = 21.232876712328
$ws.Range("F17").Value = 687
$ws.Range("G17").Value = 652
$ws.Range("H17").Value = 5.368098159509
$ws.Range("I17").Value = 6763
$ws.Range("J17").Value = 6459
$ws.Range("K17").Value = 4.706610930484
$ws.Range("L17").Value = 7.554071246819
$ws.Range("M17").Value = 95.915411355735
$ws.Range("N17").Value = -4.288140390602
$ws.Range("C18").Value = 50
$ws.Range("D18").Value = 66
$ws.Range("E18").Value = -24.242424242424
$ws.Range("F18").Value = 193
$ws.Range("G18").Value = 239
$ws.Range("H18").Value = -19.246861924686
$ws.Range("I18").Value = 2135
$ws.Range("J18").Value = 2261
$ws.Range("K18").Value = -5.572755417956
$ws.Range("L18").Value = -5.572755417956
$ws.Range("M18").Value = -15.210484511517
$ws.Range("N18").Value = -85.120914349432
$ws.Range("C19").Value = 175
$ws.Range("D19").Value = 208
$ws.Range("E19").Value = -15.865384615384
$ws.Range("F19").Value = 738
$ws.Range("G19").Value = 761
$ws.Range("H19").Value = -3.022339027595
$ws.Range("I19").Value = 7061
$ws.Range("J19").Value = 7184
$ws.Range("K19").Value = -1.712138084632
$ws.Range("L19").Value = 15.056216392374
$ws.Range("M19").Value = 98.120089786756
$ws.Range("N19").Value = 21.052631578947
$ws.Range("C20").Value = 65
$ws.Range("D20").Value = 70
$ws.Range("E20").Value = -7.142857142857
$ws.Range("F20").Value = 299
$ws.Range("G20").Value = 325
$ws.Range("H20").Value = -8
$ws.Range("I20").Value = 3418
$ws.Range("J20").Value = 3239
$ws.Range("K20").Value = 5.526397036122
$ws.Range("L20").Value = -16.060903732809
$ws.Range("M20").Value = 114.698492462312
$ws.Range("N20").Value = -70.806286299965
$ws.Range("C21").Value = 568
$ws.Range("D21").Value = 602
$ws.Range("E21").Value = -5.647840531561
$ws.Range("F21").Value = 2372
$ws.Range("G21").Value = 2427
$ws.Range("H21").Value = -2.266172229089
$ws.Range("I21").Value = 23591
$ws.Range("J21").Value = 23394
$ws.Range("K21").Value = 0.842096263999
$ws.Range("L21").Value = 3.116531165311
$ws.Range("M21").Value = 59.14058283864
$ws.Range("N21").Value = -54.757977907333
$ws.Range("C22").Value = 7
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 75
$ws.Range("F22").Value = 20
$ws.Range("G22").Value = 24
$ws.Range("H22").Value = -16.666666666666
$ws.Range("I22").Value = 232
$ws.Range("J22").Value = 259
$ws.Range("K22").Value = -10.424710424710
$ws.Range("L22").Value = 2.654867256637
$ws.Range("M22").Value = -1.694915254237
$ws.Range("C23").Value = 37
$ws.Range("D23").Value = 21
$ws.Range("E23").Value = 76.190476190476
$ws.Range("F23").Value = 136
$ws.Range("G23").Value = 115
$ws.Range("H23").Value = 18.260869565217
$ws.Range("I23").Value = 1211
$ws.Range("J23").Value = 1333
$ws.Range("K23").Value = -9.152288072018
$ws.Range("L23").Value = -11.411850768105
$ws.Range("M23").Value = 47.682926829268
$ws.Range("C24").Value = 317
$ws.Range("D24").Value = 375
$ws.Range("E24").Value = -15.466666666666
$ws.Range("F24").Value = 1365
$ws.Range("G24").Value = 1338
$ws.Range("H24").Value = 2.017937219730
$ws.Range("I24").Value = 13644
$ws.Range("J24").Value = 12442
$ws.Range("K24").Value = 9.660826233724
$ws.Range("L24").Value = -0.871839581517
$ws.Range("M24").Value = 39.723502304147
$ws.Range("C25").Value = 99
$ws.Range("D25").Value = 151
$ws.Range("E25").Value = -34.437086092715
$ws.Range("F25").Value = 403
$ws.Range("G25").Value = 508
$ws.Range("H25").Value = -20.669291338582
$ws.Range("I25").Value = 4486
$ws.Range("J25").Value = 4937
$ws.Range("K25").Value = -9.135102288839
$ws.Range("L25").Value = -23.590529722364
$ws.Range("C26").Value = 215
$ws.Range("D26").Value = 243
$ws.Range("E26").Value = -11.522633744856
$ws.Range("F26").Value = 878
$ws.Range("G26").Value = 936
$ws.Range("H26").Value = -6.196581196581
$ws.Range("I26").Value = 8576
$ws.Range("J26").Value = 8590
$ws.Range("K26").Value = -0.162980209545
$ws.Range("L26").Value = 6.959341481666
$ws.Range("M26").Value = -0.162980209545
$ws.Range("C27").Value = 15
$ws.Range("D27").Value = 8
$ws.Range("E27").Value = 87.5
$ws.Range("F27").Value = 63
$ws.Range("G27").Value = 41
$ws.Range("H27").Value = 53.658536585365
$ws.Range("I27").Value = 526
$ws.Range("J27").Value = 482
$ws.Range("K27").Value = 9.128630705394
$ws.Range("L27").Value = 4.780876494023
$ws.Range("C28").Value = 16
$ws.Range("D28").Value = 27
$ws.Range("E28").Value = -40.740740740740
$ws.Range("F28").Value = 78
$ws.Range("G28").Value = 88
$ws.Range("H28").Value = -11.363636363636
$ws.Range("I28").Value = 815
$ws.Range("J28").Value = 914
$ws.Range("K28").Value = -10.831509846827
$ws.Range("L28").Value = 0.122850122850
$ws.Range("C29").Value = 6
$ws.Range("D29").Value = 9
$ws.Range("E29").Value = -33.333333333333
$ws.Range("F29").Value = 24
$ws.Range("G29").Value = 40
$ws.Range("H29").Value = -40
$ws.Range("I29").Value = 263
$ws.Range("J29").Value = 337
$ws.Range("K29").Value = -21.958456973293
$ws.Range("L29").Value = -17.034700315457
$ws.Range("M29").Value = -30.606860158311
$ws.Range("N29").Value = -76.263537906137
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 8
$ws.Range("E30").Value = -62.5
$ws.Range("F30").Value = 17
$ws.Range("G30").Value = 32
$ws.Range("H30").Value = -46.875
$ws.Range("I30").Value = 210
$ws.Range("J30").Value = 265
$ws.Range("K30").Value = -20.754716981132
$ws.Range("L30").Value = -19.540229885057
$ws.Range("M30").Value = -33.962264150943
$ws.Range("N30").Value = -78.809283551967
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 17
$ws.Range("J31").Value = 19
$ws.Range("K31").Value = -10.526315789473
$ws.Range("L31").Value = 6.25
$ws.Range("H33").Value = -100
$ws.Range("I33").Value = 23
$ws.Range("J33").Value = 40
$ws.Range("K33").Value = -42.5
$ws.Range("L33").Value = -37.837837837837

# --- Special cells (rows 31 and 33) that change data type (text <-> number) ---
# These require copying both the number format/style and the value/type together,
# since directly assigning .Value would either coerce text to a number or create a
# brand-new (unwanted) style/number-format entry.

# C31: was the text placeholder "0" -> becomes numeric 2 (use D14's numeric style)
$ws.Range("D14").Copy() | Out-Null
$ws.Range("C31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("C31").Value = 2

# F31: was the text placeholder "0" -> becomes numeric 2 (use D14's numeric style)
$ws.Range("D14").Copy() | Out-Null
$ws.Range("F31").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("F31").Value = 2

# F33: was numeric 3 -> becomes the text placeholder "0" (use D33's text style + value)
$ws.Range("D33").Copy() | Out-Null
$ws.Range("F33").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D33").Copy() | Out-Null
$ws.Range("F33").PasteSpecial(-4163) | Out-Null   # xlPasteValues (copies text type + value)

# G33: was numeric 5 -> becomes numeric 3 (style unchanged)
$ws.Range("G33").Value = 3

$excel.CutCopyMode = 0
